# Generate Report for Handoff
# Adds a new handed-off file (b095d723-d7d4-4dc8-96cf-d41fd8da1740.md) as
# row 9 on all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileGuid = "b095d723-d7d4-4dc8-96cf-d41fd8da1740"
$fileName = "$fileGuid.md"
$pathAndName = "e2e\$fileName"
$commitSha = "6bfb1e39b49dcd2ab32c5c5e6f422e85ba37ccd2"
$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob"
$newFileSha = "a9c6e6b6cf7b6a6e2c1b4a1e7f9d2c3b5a6e7f8a"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")
$loO = $wsO.ListObjects.Item(1)
$loO.ListRows.Add() | Out-Null

$wsO.Range("A9").Value = $fileName
$wsO.Hyperlinks.Add($wsO.Range("B9"), "$repoBase/$newFileSha/$pathAndName", "", "", $pathAndName)
$wsO.Range("C9").Value = ".md"
$wsO.Range("D9").Value = ""
$wsO.Range("E9").Value = "Ready for handoff"
$wsO.Range("F9").Value = "Ready for handoff"
$wsO.Range("G9").Value = "2016-09-03 18:56:45"
$wsO.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")
$loZ = $wsZ.ListObjects.Item(1)
$loZ.ListRows.Add() | Out-Null

$wsZ.Hyperlinks.Add($wsZ.Range("A9"), "$repoBase/$newFileSha/$pathAndName", "", "", $fileName)
$wsZ.Range("B9").Value = ".md"
$wsZ.Range("C9").Value = "Ready for handoff"
$wsZ.Range("D9").Value = "e2e"
$wsZ.Range("E9").Value = "ht"
$wsZ.Range("F9").Value = "False"
$wsZ.Range("G9").Value = "$fileGuid.$commitSha.zh-cn.xlf"
$wsZ.Range("H9").Value = "2016-09-03 18:56:41"
$wsZ.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("I9").Value = ""
$wsZ.Range("J9").Value = ""
$wsZ.Range("K9").Value = "0001-01-01 00:00:00"
$wsZ.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("L9").Value = ""
$wsZ.Range("M9").Value = "True"
$wsZ.Range("N9").Value = ""
$wsZ.Range("O9").Value = "False"
$wsZ.Range("P9").Value = ""

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")
$loD = $wsD.ListObjects.Item(1)
$loD.ListRows.Add() | Out-Null

$wsD.Hyperlinks.Add($wsD.Range("A9"), "$repoBase/$newFileSha/$pathAndName", "", "", $fileName)
$wsD.Range("B9").Value = ".md"
$wsD.Range("C9").Value = "Ready for handoff"
$wsD.Range("D9").Value = "e2e"
$wsD.Range("E9").Value = "ht"
$wsD.Range("F9").Value = "False"
$wsD.Range("G9").Value = "$fileGuid.$commitSha.de-de.xlf"
$wsD.Range("H9").Value = "2016-09-03 18:56:45"
$wsD.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("I9").Value = ""
$wsD.Range("J9").Value = ""
$wsD.Range("K9").Value = "0001-01-01 00:00:00"
$wsD.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("L9").Value = ""
$wsD.Range("M9").Value = "True"
$wsD.Range("N9").Value = ""
$wsD.Range("O9").Value = "False"
$wsD.Range("P9").Value = ""

Write-Host "Row 9 added to Overview, zh-cn, de-de sheets."
